$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Data" -> "Summary"
$ws.Name = "Summary"

# Insert 4 new blank rows before the existing row 5, shifting the table
# (old rows 5-9) down to rows 9-13, opening up space for the new
# "Source Type" line at row 7.
$ws.Rows("5:8").Insert()

# --- Re-assert formatting on the cells that were already in the sheet and
# just shifted down a few rows. (Re-setting the font explicitly forces a
# clean font/style entry instead of a stale/shared one.) ---
$ws.Range("A1").Font.Size = 18

$ws.Range("A3").Font.Bold = $true

$ws.Range("B9:D9").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("A11").Font.Bold = $true
$ws.Range("A12").Font.Bold = $true

$ws.Range("A13").Font.Italic = $true

# New "Source Type" line in the freshly-inserted row 7, bold + underlined
# (the new "title_" style).
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true
$ws.Range("A7").Font.Italic = $false

# New source-detail lines appended at the bottom of the sheet.
$ws.Range("A21").Value = "NBS"
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").Font.Underline = $false
$ws.Range("A21").Font.Italic = $false

$ws.Range("A22").Value = "National Bureau of Statistics (NBS), ""Statistical Yearbook for Southern Sudan"", 2010, p. 142. Available at http://ssnbs.org/statistical-year-book/"
$ws.Range("A22").Font.Italic = $true
$ws.Range("A22").Font.Bold = $false
$ws.Range("A22").Font.Underline = $false
